# "Updated: po 23. 08. 2021" -- refresh AgTests/AgPosit (cols F/G) corrections
# across the historical series, and append three new daily rows
# (2021-08-20, 2021-08-21, 2021-08-22).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Corrections to existing AgTests (F) / AgPosit (G) values ---
$ws.Range("F331").Value = 154371
$ws.Range("F334").Value = 192925
$ws.Range("F336").Value = 82023
$ws.Range("F337").Value = 105633
$ws.Range("F342").Value = 178925
$ws.Range("F343").Value = 134000
$ws.Range("F377").Value = 177043
$ws.Range("F380").Value = 346014
$ws.Range("F384").Value = 172541
$ws.Range("F385").Value = 151274
$ws.Range("F386").Value = 183527
$ws.Range("F391").Value = 178272
$ws.Range("F398").Value = 300871
$ws.Range("G399").Value = 968
$ws.Range("F400").Value = 150053
$ws.Range("F405").Value = 175109
$ws.Range("F406").Value = 171786
$ws.Range("F407").Value = 158563
$ws.Range("F408").Value = 306114
$ws.Range("F411").Value = 225644
$ws.Range("F412").Value = 177049
$ws.Range("F413").Value = 149969
$ws.Range("F418").Value = 202735
$ws.Range("F419").Value = 149736
$ws.Range("F420").Value = 139370
$ws.Range("F421").Value = 153633
$ws.Range("F422").Value = 299017
$ws.Range("F425").Value = 138314
$ws.Range("F427").Value = 90695
$ws.Range("F428").Value = 102736
$ws.Range("F429").Value = 178967
$ws.Range("F433").Value = 87466
$ws.Range("F434").Value = 79497
$ws.Range("F435").Value = 83682
$ws.Range("F436").Value = 145557
$ws.Range("F439").Value = 89496
$ws.Range("F440").Value = 73931
$ws.Range("F441").Value = 68550
$ws.Range("F442").Value = 70769
$ws.Range("F453").Value = 70392
$ws.Range("F454").Value = 52769
$ws.Range("F470").Value = 43642
$ws.Range("F501").Value = 5783
$ws.Range("F502").Value = 10658
$ws.Range("F503").Value = 7594
$ws.Range("F504").Value = 7578
$ws.Range("F505").Value = 8650
$ws.Range("F506").Value = 11019
$ws.Range("F507").Value = 7338
$ws.Range("F508").Value = 5797
$ws.Range("F509").Value = 9729
$ws.Range("F510").Value = 7948
$ws.Range("F511").Value = 6897
$ws.Range("F512").Value = 8589
$ws.Range("F513").Value = 10480
$ws.Range("F514").Value = 7095
$ws.Range("F515").Value = 5138
$ws.Range("F516").Value = 9451
$ws.Range("F517").Value = 6824
$ws.Range("F518").Value = 7177
$ws.Range("G518").Value = 10
$ws.Range("F519").Value = 7975
$ws.Range("F520").Value = 10319
$ws.Range("F521").Value = 6806
$ws.Range("F522").Value = 5084
$ws.Range("F523").Value = 10170
$ws.Range("F524").Value = 7809
$ws.Range("F525").Value = 7597
$ws.Range("F526").Value = 8739
$ws.Range("F527").Value = 11380
$ws.Range("F528").Value = 7981
$ws.Range("G528").Value = 21
$ws.Range("F529").Value = 5572
$ws.Range("F530").Value = 12417
$ws.Range("G530").Value = 40
$ws.Range("F531").Value = 8754
$ws.Range("G531").Value = 25
$ws.Range("F532").Value = 9643
$ws.Range("G532").Value = 51
$ws.Range("F533").Value = 11028
$ws.Range("G533").Value = 38

# --- Append three new daily rows (534-536) ---
$ws.Range("A534").Value = 44428
$ws.Range("B534").Value = 393977
$ws.Range("C534").Value = 7484
$ws.Range("D534").Value = 94
$ws.Range("E534").Value = 12547
$ws.Range("F534").Value = 14171
$ws.Range("G534").Value = 47

$ws.Range("A535").Value = 44429
$ws.Range("B535").Value = 394082
$ws.Range("C535").Value = 4842
$ws.Range("D535").Value = 105
$ws.Range("E535").Value = 12547
$ws.Range("F535").Value = 7952
$ws.Range("G535").Value = 19

$ws.Range("A536").Value = 44430
$ws.Range("B536").Value = 394093
$ws.Range("C536").Value = 943
$ws.Range("D536").Value = 11
$ws.Range("E536").Value = 12547
$ws.Range("F536").Value = 5262
$ws.Range("G536").Value = 31
